# Add_Supplier_Afghanistan.xlsx regression-data refresh
# ("supplierregress0N" replacing the stale "Dailyregress0N" test accounts,
#  plus the CoC request URL / invite-status refresh and a changed selection.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add_Supplier sheet: rename the regression accounts used in columns A/B/C
#    and flip the "User_Invited"/"User_Registered" flags in row 2/3 to "No".
# ---------------------------------------------------------------------------
$wsSupplier = $wb.Worksheets.Item("Add_Supplier")

$accounts = @(
    @{ Row = 2; User = "supplierregress01" },
    @{ Row = 3; User = "supplierregress02" },
    @{ Row = 4; User = "supplierregress03" },
    @{ Row = 5; User = "supplierregress04" },
    @{ Row = 6; User = "supplierregress05" },
    @{ Row = 7; User = "supplierregress06" },
    @{ Row = 8; User = "supplierregress07" }
)

foreach ($acct in $accounts) {
    $row = $acct.Row
    $user = $acct.User
    $email = "$user@yopmail.com"

    $wsSupplier.Range("A$row").Value = $user
    $wsSupplier.Range("B$row").Value = $email
    $wsSupplier.Range("C$row").Value = $user
}

# The invite/registration indicator columns for the first two accounts move
# from "Yes" to "No".
$wsSupplier.Range("G2").Value = "No"
$wsSupplier.Range("H2").Value = "No"
$wsSupplier.Range("G3").Value = "No"

# Keep each mailto hyperlink's visible text in sync with its own cell, and
# repoint every ScreenTip at the refreshed "master" account (row 2).
foreach ($link in $wsSupplier.Hyperlinks) {
    $addr = $link.Range().Address()
    $addr = $addr.Replace("$", "")
    $cellValue = $wsSupplier.Range($addr).Value()
    $link.TextToDisplay = $cellValue
    $link.ScreenTip = "mailto:supplierregress01@yopmail.com"
}

# The saved selection moves from H2:H8 to a single A2.
$wsSupplier.Activate()
$wsSupplier.Range("A2").Select()

# ---------------------------------------------------------------------------
# 2) Recent_Added sheet: refresh the stale hyperlink display text so it
#    matches the account actually shown in B2.
# ---------------------------------------------------------------------------
$wsRecent = $wb.Worksheets.Item("Recent_Added")
foreach ($link in $wsRecent.Hyperlinks) {
    $addr = $link.Range().Address()
    $addr = $addr.Replace("$", "")
    $cellValue = $wsRecent.Range($addr).Value()
    $link.TextToDisplay = $cellValue
}
